$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "766×5=3830" "187×5=935"
Replace-Text "606×7=4242" "109×2=218"
Replace-Text "894×8=7152" "795×3=2385"
Replace-Text "531×8=4248" "145×3=435"
Replace-Text "603×8=4824" "226×7=1582"
Replace-Text "111×9=999" "132×6=792"
Replace-Text "889×2=1778" "259×4=1036"
Replace-Text "286×9=2574" "406×6=2436"
Replace-Text "872×2=1744" "592×4=2368"
Replace-Text "514×7=3598" "234×7=1638"
Replace-Text "219×8=1752" "139×4=556"
Replace-Text "493×3=1479" "938×7=6566"
Replace-Text "310×9=2790" "113×7=791"
Replace-Text "317×8=2536" "834×8=6672"
Replace-Text "832×3=2496" "481×6=2886"
Replace-Text "362×5=1810" "785×9=7065"
Replace-Text "420×2=840" "658×8=5264"
Replace-Text "983×8=7864" "406×8=3248"
Replace-Text "944×6=5664" "668×4=2672"
Replace-Text "785×3=2355" "263×5=1315"
Replace-Text "844×7=5908" "893×7=6251"
Replace-Text "323×3=969" "812×4=3248"
Replace-Text "914×8=7312" "351×7=2457"
Replace-Text "683×8=5464" "339×2=678"
Replace-Text "685×9=6165" "785×2=1570"
